# WebScanPro Final PPT.pptx - "Add files via upload" edit
#
# Summary of the applied change (see xml_diff):
#   - Slide 1 title box ("WebScanPro: AI-Powered Vulnerability Assessment"):
#       widened slightly and the text's font size nudged down 44.5pt -> 43.5pt.
#   - A handful of autoshapes across slides 3, 4, 6, 8 and 10 were
#       repositioned / resized by small amounts (sub-point rounding tweaks
#       plus a handful of deliberate width increases, e.g. text boxes that
#       got wider so their text would not wrap).
#
# NOTE on numeric literals below: PowerPoint's Shape.Left/Top/Width/Height
# are expressed in points (1 pt = 12700 EMU) and are stored by the host as
# single-precision floats, then converted back to EMU by truncation. The
# literals here were chosen so that, after that float32 round-trip, they
# land exactly on the target EMU value from the target OOXML.

$p = $ppt.ActivePresentation

# --- Slide 1: title textbox - widen + shrink font from 44.5pt to 43.5pt ---
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Google Shape;56;p13")
$sh.Width = 1084.7244873046875

$tr = $sh.TextFrame.TextRange
$runCount = $tr.Runs().Count
for ($i = 1; $i -le $runCount; $i++) {
    $tr.Runs($i).Font.Size = 43.5
}

# --- Slide 3: "Target Scanning Module Development" header ---
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item("Google Shape;101;p15")
$sh.Left = 35.58622360229492
$sh.Top = 59.9527587890625
$sh.Width = 618.3307495117188
$sh.Height = 31.77165412902832

# --- Slide 4: "DVWA Exploitation Results" header ---
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item("Google Shape;162;p16")
$sh.Left = 19.560789108276367
$sh.Top = 268.1437072753906
$sh.Width = 211.96063232421875

# --- Slide 6: "Severity Classification" header ---
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item("Google Shape;240;p18")
$sh.Left = 50.78314971923828
$sh.Top = 321.7224426269531
$sh.Width = 281.905517578125
$sh.Height = 27.21259880065918

# --- Slide 8: "Result Aggregation" label ---
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item("Google Shape;314;p20")
$sh.Left = -16.10480499267578
$sh.Top = 121.89567565917969
$sh.Width = 203.03150939941406

# --- Slide 8: "Severity Prediction" label ---
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item("Google Shape;331;p20")
$sh.Left = 47.119686126708984
$sh.Top = 408.4626159667969
$sh.Width = 162.44882202148438

# --- Slide 8: empty rounded-rect card, nudged down ---
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item("Google Shape;340;p20")
$sh.Top = 332.5731506347656

# --- Slide 10: small status-dot marker ---
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item("Google Shape;423;p22")
$sh.Top = 228.06552124023438
$sh.Width = 5.952755928039551
$sh.Height = 5.952755928039551

# --- Slide 10: "AI Integration" label ---
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item("Google Shape;424;p22")
$sh.Top = 221.74205017089844
$sh.Width = 148.93701171875
$sh.Height = 18.614173889160156

# --- Slide 10: "Machine learning-powered payload generation..." body text ---
$s = $p.Slides.Item(10)
$sh = $s.Shapes.Item("Google Shape;425;p22")
$sh.Top = 252.26708984375
$sh.Width = 458.0079040527344
$sh.Height = 19.0629940032959
